$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - match style of existing header cells (e.g. E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Data cells F2:F16 - text values (time_taken timestamps)
$timeValues = @(
    "2021-10-05 13:38:42.714838",
    "2021-10-05 13:38:42.714851",
    "2021-10-05 13:38:42.714855",
    "2021-10-05 13:38:42.714858",
    "2021-10-05 13:38:42.714862",
    "2021-10-05 13:38:42.714865",
    "2021-10-05 13:38:42.714868",
    "2021-10-05 13:38:42.714871",
    "2021-10-05 13:38:42.714874",
    "2021-10-05 13:38:42.714877",
    "2021-10-05 13:38:42.714880",
    "2021-10-05 13:38:42.714883",
    "2021-10-05 13:38:42.714886",
    "2021-10-05 13:38:42.714889",
    "2021-10-05 13:38:42.714892"
)

for ($i = 0; $i -lt $timeValues.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $timeValues[$i]
}
